$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 867.2857
$ws.Range("I12").Value = 774.2
$ws.Range("J12").Value = 1100
$ws.Range("K12").Value = 774.2
$ws.Range("L12").Value = 1100
$ws.Range("M12").Value = -604.2
$ws.Range("N12").Value = -1440
$ws.Range("H43").Value = 1744.4762
$ws.Range("I43").Value = 1413.9231
$ws.Range("J43").Value = 2281.625
$ws.Range("K43").Value = 1413.9231
$ws.Range("L43").Value = 2281.625
$ws.Range("M43").Value = -1344.9231
$ws.Range("N43").Value = -2419.625
$ws.Range("H86").Value = 4815.5
$ws.Range("I86").Value = 4149
$ws.Range("J86").Value = 5325.1763
$ws.Range("K86").Value = 4149
$ws.Range("L86").Value = 5325.1763
$ws.Range("M86").Value = -3026
$ws.Range("N86").Value = -7571.1763
$ws.Range("H89").Value = 4815.5
$ws.Range("I89").Value = 4149
$ws.Range("J89").Value = 5325.1763
$ws.Range("K89").Value = 20745
$ws.Range("L89").Value = 26625.8815
$ws.Range("M89").Value = -15129
$ws.Range("N89").Value = -37857.8815
$ws.Range("H107").Value = 55556204
$ws.Range("I107").Value = 111111250
$ws.Range("J107").Value = 1166.3334
$ws.Range("K107").Value = 111111250
$ws.Range("L107").Value = 1166.3334
$ws.Range("M107").Value = -111109330
$ws.Range("N107").Value = -5006.3334
$ws.Range("H137").Value = 130808.71
$ws.Range("I137").Value = 180032.2
$ws.Range("J137").Value = 7750
$ws.Range("K137").Value = 540096.6000000001
$ws.Range("L137").Value = 23250
$ws.Range("M137").Value = -537546.6000000001
$ws.Range("N137").Value = -28350
$ws.Range("H138").Value = 5774.7188
$ws.Range("J138").Value = 5992.9644
$ws.Range("L138").Value = 17978.8932
$ws.Range("N138").Value = -28258.8932
$ws.Range("H141").Value = 6346.593
$ws.Range("I141").Value = 6682.2915
$ws.Range("K141").Value = 20046.8745
$ws.Range("M141").Value = -14866.8745

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10207768
$ws.Range("I45").Value = 17858354
$ws.Range("K45").Value = 17858354
$ws.Range("M45").Value = -17857977
$ws.Range("H61").Value = 4653.1787
$ws.Range("I61").Value = 4640.5186
$ws.Range("K61").Value = 4640.5186
$ws.Range("M61").Value = -4428.5186
$ws.Range("H63").Value = 8100.727
$ws.Range("I63").Value = 2602
$ws.Range("J63").Value = 9322.666999999999
$ws.Range("K63").Value = 2602
$ws.Range("L63").Value = 9322.666999999999
$ws.Range("M63").Value = -1916
$ws.Range("N63").Value = -10694.667
$ws.Range("H66").Value = 8100.727
$ws.Range("I66").Value = 2602
$ws.Range("J66").Value = 9322.666999999999
$ws.Range("K66").Value = 13010
$ws.Range("L66").Value = 46613.335
$ws.Range("M66").Value = -9578
$ws.Range("N66").Value = -53477.335
$ws.Range("H110").Value = 1744516.8
$ws.Range("I110").Value = 2138358.2
$ws.Range("K110").Value = 2138358.2
$ws.Range("M110").Value = -2136313.2
$ws.Range("H122").Value = 11831191
$ws.Range("I122").Value = 19633926
$ws.Range("K122").Value = 58901778
$ws.Range("M122").Value = -58899328
$ws.Range("H132").Value = 30489.861
$ws.Range("I132").Value = 2644.5881
$ws.Range("J132").Value = 55404.05
$ws.Range("K132").Value = 7933.7643
$ws.Range("L132").Value = 166212.15
$ws.Range("M132").Value = -5403.7643
$ws.Range("N132").Value = -171272.15
$ws.Range("H136").Value = 4653.1787
$ws.Range("I136").Value = 4640.5186
$ws.Range("K136").Value = 13921.5558
$ws.Range("M136").Value = -11371.5558

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1882.6316
$ws.Range("I20").Value = 1979.7142
$ws.Range("J20").Value = 1762.7059
$ws.Range("K20").Value = 1979.7142
$ws.Range("L20").Value = 1762.7059
$ws.Range("M20").Value = -1732.7142
$ws.Range("N20").Value = -2256.7059
$ws.Range("H82").Value = 6166.6665
$ws.Range("I82").Value = 2400
$ws.Range("K82").Value = 2400
$ws.Range("M82").Value = -2017
$ws.Range("H85").Value = 6166.6665
$ws.Range("I85").Value = 2400
$ws.Range("K85").Value = 2400
$ws.Range("M85").Value = -1074

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6633.4
$ws.Range("I4").Value = 5001
$ws.Range("J4").Value = 7041.5
$ws.Range("K4").Value = 5001
$ws.Range("L4").Value = 7041.5
$ws.Range("M4").Value = -4889
$ws.Range("N4").Value = -7265.5
$ws.Range("H31").Value = 47266.047
$ws.Range("I31").Value = 4585.3335
$ws.Range("J31").Value = 76814.234
$ws.Range("K31").Value = 4585.3335
$ws.Range("L31").Value = 76814.234
$ws.Range("M31").Value = -4290.3335
$ws.Range("N31").Value = -77404.234
$ws.Range("H34").Value = 47266.047
$ws.Range("I34").Value = 4585.3335
$ws.Range("J34").Value = 76814.234
$ws.Range("K34").Value = 4585.3335
$ws.Range("L34").Value = 76814.234
$ws.Range("M34").Value = -4383.3335
$ws.Range("N34").Value = -77218.234
$ws.Range("H58").Value = 6996.206
$ws.Range("I58").Value = 8333.458000000001
$ws.Range("K58").Value = 8333.458000000001
$ws.Range("M58").Value = -8130.458000000001
$ws.Range("H132").Value = 72390.14999999999
$ws.Range("I132").Value = 68589.2
$ws.Range("J132").Value = 77141.336
$ws.Range("K132").Value = 205767.6
$ws.Range("L132").Value = 231424.008
$ws.Range("M132").Value = -203237.6
$ws.Range("N132").Value = -236484.008
$ws.Range("H136").Value = 6996.206
$ws.Range("I136").Value = 8333.458000000001
$ws.Range("K136").Value = 25000.374
$ws.Range("M136").Value = -22450.374
$ws.Range("H141").Value = 558000.3
$ws.Range("J141").Value = 558000.3
$ws.Range("L141").Value = 558000.3
$ws.Range("N141").Value = -568360.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H133").Value = 8113.778
$ws.Range("I133").Value = 3030
$ws.Range("J133").Value = 8749.25
$ws.Range("K133").Value = 9090
$ws.Range("L133").Value = 26247.75
$ws.Range("M133").Value = -4030
$ws.Range("N133").Value = -36367.75
$ws.Range("H137").Value = 3101.889
$ws.Range("I137").Value = 2319.6667
$ws.Range("J137").Value = 4666.3335
$ws.Range("K137").Value = 6959.000100000001
$ws.Range("L137").Value = 13999.0005
$ws.Range("M137").Value = -1859.000100000001
$ws.Range("N137").Value = -24199.0005

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18195502
$ws.Range("I70").Value = 22226068
$ws.Range("K70").Value = 22226068
$ws.Range("M70").Value = -22225798
$ws.Range("H73").Value = 18195502
$ws.Range("I73").Value = 22226068
$ws.Range("K73").Value = 22226068
$ws.Range("M73").Value = -22225132
$ws.Range("H80").Value = 1842835.2
$ws.Range("I80").Value = 2629672.5
$ws.Range("K80").Value = 2629672.5
$ws.Range("M80").Value = -2628674.5
$ws.Range("H83").Value = 1842835.2
$ws.Range("I83").Value = 2629672.5
$ws.Range("K83").Value = 13148362.5
$ws.Range("M83").Value = -13143370.5
$ws.Range("H97").Value = 882556.7
$ws.Range("I97").Value = 1587935.1
$ws.Range("K97").Value = 1587935.1
$ws.Range("M97").Value = -1587439.1
$ws.Range("H122").Value = 164313.42
$ws.Range("I122").Value = 195555.38
$ws.Range("J122").Value = 4632.3335
$ws.Range("K122").Value = 586666.14
$ws.Range("L122").Value = 13897.0005
$ws.Range("M122").Value = -584216.14
$ws.Range("N122").Value = -18797.0005
$ws.Range("H132").Value = 4936.952
$ws.Range("I132").Value = 3977.7144
$ws.Range("J132").Value = 6855.4287
$ws.Range("K132").Value = 11933.1432
$ws.Range("L132").Value = 20566.2861
$ws.Range("M132").Value = -9403.143199999999
$ws.Range("N132").Value = -25626.2861

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9192
$ws.Range("I7").Value = 8232.4
$ws.Range("K7").Value = 8232.4
$ws.Range("M7").Value = -8120.4
$ws.Range("H16").Value = 219.89655
$ws.Range("I16").Value = 168.76923
$ws.Range("K16").Value = 168.76923
$ws.Range("M16").Value = 1.230770000000007
$ws.Range("H22").Value = 82925.91
$ws.Range("I22").Value = 149614.33
$ws.Range("K22").Value = 149614.33
$ws.Range("M22").Value = -149319.33
$ws.Range("H27").Value = 82925.91
$ws.Range("I27").Value = 149614.33
$ws.Range("K27").Value = 149614.33
$ws.Range("M27").Value = -149507.33
$ws.Range("H68").Value = 3758.4285
$ws.Range("I68").Value = 2100.6667
$ws.Range("J68").Value = 5001.75
$ws.Range("K68").Value = 2100.6667
$ws.Range("L68").Value = 5001.75
$ws.Range("M68").Value = -1351.6667
$ws.Range("N68").Value = -6499.75
$ws.Range("H71").Value = 3758.4285
$ws.Range("I71").Value = 2100.6667
$ws.Range("J71").Value = 5001.75
$ws.Range("K71").Value = 10503.3335
$ws.Range("L71").Value = 25008.75
$ws.Range("M71").Value = -6759.333500000001
$ws.Range("N71").Value = -32496.75
$ws.Range("H126").Value = 9192
$ws.Range("I126").Value = 8232.4
$ws.Range("K126").Value = 24697.2
$ws.Range("M126").Value = -22227.2
$ws.Range("H132").Value = 11285.071
$ws.Range("I132").Value = 12047.571
$ws.Range("J132").Value = 8997.571
$ws.Range("K132").Value = 36142.713
$ws.Range("L132").Value = 26992.713
$ws.Range("M132").Value = -33612.713
$ws.Range("N132").Value = -32052.713

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 30304142
$ws.Range("I107").Value = 40000590
$ws.Range("K107").Value = 120001770
$ws.Range("M107").Value = -119999850
$ws.Range("H136").Value = 3680.244
$ws.Range("I136").Value = 3323.7742
$ws.Range("K136").Value = 9971.3226
$ws.Range("M136").Value = -7421.3226
